# Oppdatert storre naering for levering til nettet ++
#
# The "Storre naeringskunde" (bigger business customer) sheet gets a new
# price row: the old single "Energiledd" line (summer+winter lumped
# together) is split into two separate lines - "Energiledd sommer" and
# "Energiledd vinter" - each with its own (higher) price. Also the
# previously-selected sheet/tab moves back to the first sheet.

$wb = $excel.ActiveWorkbook

$wsPrivat = $wb.Worksheets.Item(1)   # Privatkunde
$wsStorre = $wb.Worksheets.Item(3)   # Storre naeringskunde

# --- Storre naeringskunde: split "Energiledd" into summer/winter rows ---

# Insert a fresh row above the current "Energiledd" row (row 4). Excel
# shifts everything below down by one and keeps formulas/merges consistent.
$wsStorre.Rows.Item(4).Insert()

# The new blank row 4 should look like the "Forbruksavgift" rows further
# down (plain text style, no border on the unit column for the summer
# line) - copy that formatting across before filling in the values.
$wsStorre.Range("A10:D10").Copy()
$wsStorre.Range("A4:D4").PasteSpecial(-4122)

# New row 4: "Energiledd sommer"
$wsStorre.Cells.Item(4, 1).Value2 = "Energiledd sommer"
$wsStorre.Cells.Item(4, 2).Value2 = 7.3
$wsStorre.Cells.Item(4, 3).Value2 = 4.0999999999999996
$wsStorre.Cells.Item(4, 4).ClearContents()

# Row 5 (the old "Energiledd" row, pushed down by the insert) becomes
# "Energiledd vinter" with the same new prices; it keeps its own existing
# formatting (unit column still shows "ore/kWh").
$wsStorre.Cells.Item(5, 1).Value2 = "Energiledd vinter"
$wsStorre.Cells.Item(5, 2).Value2 = 7.3
$wsStorre.Cells.Item(5, 3).Value2 = 4.0999999999999996

# Restore the selection on this sheet to B7 (where the old row 8, now at
# row 7, keeps getting edited).
[void]$wsStorre.Range("B7").Select()

# --- Active tab moves back to the first sheet ---
[void]$wsPrivat.Activate()
